$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("translation")

$ws.Range("A9").Value = "Enter One to Ones"
$ws.Range("B9").Value = "Enter One to Ones"
$ws.Range("C9").Value = "Enter One to Ones"
$ws.Range("D9").Value = "Enter One to Ones"
